# The deck originally ships with two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (the stock blue/orange Office palette)
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet" (the palette actually used
#                            by the slide master / whole deck design)
#
# The authored change swaps the two themes' contents: the design that drives
# the slides/slide master switches from the "Integral" (Red Violet) palette to
# the plain "Office Theme" palette. (fontScheme/fmtScheme are identical between
# the two theme parts already, so the only observable difference is the set of
# twelve scheme colours.)
#
# PowerPoint's object model exposes the live design's colours through
# SlideMaster.Theme.ThemeColorScheme (Item 1..12 map to dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink in that order). Item(n).RGB takes/returns an OLE
# color (0xBBGGRR) the same way the classic VBA RGB() function packs it, so we
# convert each target hex colour before assigning it.

function ConvertTo-OleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Target palette: the stock "Office Theme" colours (previously theme1.xml),
# now applied to the design's live theme (theme2.xml) per the swap.
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = ConvertTo-OleColor $officeColors[$i - 1]
}
